$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Friday, Jan 13"
$ws.Range("C19").Value = "9:35 AM"
$ws.Range("D19").Value = "W92066"
$ws.Range("E19").Value = "London"
$ws.Range("F19").Value = "(LTN)"
$ws.Range("G19").Value = "Wizz Air "
$ws.Range("H19").Value = "A321"
$ws.Range("I19").Value = "(G-WUKJ)"
$ws.Range("J19").Value = "9:02 AM"
$ws.Range("K19").Interior.Pattern = -4142
$ws.Range("L19").Value = "0 hours, -33 minutes"
$ws.Range("M19").Interior.Pattern = -4142

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "Friday, Jan 13"
$ws.Range("C20").Value = "1:05 PM"
$ws.Range("D20").Value = "FR6639"
$ws.Range("E20").Value = "London"
$ws.Range("F20").Value = "(LTN)"
$ws.Range("G20").Value = "Ryanair "
$ws.Range("H20").Value = "B738"
$ws.Range("I20").Value = "(EI-DPL)"
$ws.Range("J20").Value = "12:36 PM"
$ws.Range("K20").Interior.Pattern = -4142
$ws.Range("L20").Value = "0 hours, -29 minutes"
$ws.Range("M20").Interior.Pattern = -4142
